$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: update/clear values (meanEMG legmaxROM data update)
$ws.Range("B2").Value = 9.3300853526987098
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 13.962524936529718
$ws.Range("E2").ClearContents()

# Row 3: update values
$ws.Range("B3").Value = 8.2072100330232018
$ws.Range("C3").Value = -5.4378866419480616
$ws.Range("D3").Value = 15.608761011856245
$ws.Range("E3").Value = -2.7961996361685308

# Update the selection to reflect the new active range
$ws.Range("B1:E3").Select()
